$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2501724.5
$ws.Range("I137").Value = 4763435
$ws.Range("J137").Value = 1938.9474
$ws.Range("K137").Value = 14290305
$ws.Range("L137").Value = 5816.8422
$ws.Range("M137").Value = -14287755
$ws.Range("N137").Value = -10916.8422
$ws.Range("H138").Value = 1722209.4
$ws.Range("I138").Value = 2332.8
$ws.Range("J138").Value = 2168930.8
$ws.Range("K138").Value = 6998.400000000001
$ws.Range("L138").Value = 6506792.399999999
$ws.Range("M138").Value = -1858.400000000001
$ws.Range("N138").Value = -6517072.399999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 23857898
$ws.Range("I61").Value = 33367396
$ws.Range("J61").Value = 84152.336
$ws.Range("K61").Value = 33367396
$ws.Range("L61").Value = 84152.336
$ws.Range("M61").Value = -33367184
$ws.Range("N61").Value = -84576.336
$ws.Range("H74").Value = 6462424.5
$ws.Range("I74").Value = 10040648
$ws.Range("K74").Value = 10040648
$ws.Range("M74").Value = -10039774
$ws.Range("H77").Value = 6462424.5
$ws.Range("I77").Value = 10040648
$ws.Range("K77").Value = 50203240
$ws.Range("M77").Value = -50198872
$ws.Range("H110").Value = 5005005.5
$ws.Range("I110").Value = 5005005.5
$ws.Range("K110").Value = 5005005.5
$ws.Range("M110").Value = -5002960.5
$ws.Range("H111").Value = 60644
$ws.Range("J111").Value = 60644
$ws.Range("L111").Value = 60644
$ws.Range("N111").Value = -68824
$ws.Range("H136").Value = 23857898
$ws.Range("I136").Value = 33367396
$ws.Range("J136").Value = 84152.336
$ws.Range("K136").Value = 100102188
$ws.Range("L136").Value = 252457.008
$ws.Range("M136").Value = -100099638
$ws.Range("N136").Value = -257557.008
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2641.614
$ws.Range("I134").Value = 2415.9048
$ws.Range("K134").Value = 7247.714399999999
$ws.Range("M134").Value = -4712.714399999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1385.7561
$ws.Range("I31").Value = 967.5714
$ws.Range("J31").Value = 1529.7213
$ws.Range("K31").Value = 967.5714
$ws.Range("L31").Value = 1529.7213
$ws.Range("M31").Value = -672.5714
$ws.Range("N31").Value = -2119.7213
$ws.Range("H34").Value = 1385.7561
$ws.Range("I34").Value = 967.5714
$ws.Range("J34").Value = 1529.7213
$ws.Range("K34").Value = 967.5714
$ws.Range("L34").Value = 1529.7213
$ws.Range("M34").Value = -765.5714
$ws.Range("N34").Value = -1933.7213
$ws.Range("H58").Value = 24288268
$ws.Range("I58").Value = 27570226
$ws.Range("J58").Value = 1779.6
$ws.Range("K58").Value = 27570226
$ws.Range("L58").Value = 1779.6
$ws.Range("M58").Value = -27570023
$ws.Range("N58").Value = -2185.6
$ws.Range("H132").Value = 39620.11
$ws.Range("I132").Value = 2416.4285
$ws.Range("J132").Value = 169833
$ws.Range("K132").Value = 7249.2855
$ws.Range("L132").Value = 509499
$ws.Range("M132").Value = -4719.2855
$ws.Range("N132").Value = -514559
$ws.Range("H134").Value = 59558.21
$ws.Range("I134").Value = 2642.1667
$ws.Range("K134").Value = 7926.500100000001
$ws.Range("M134").Value = -5391.500100000001
$ws.Range("H136").Value = 24288268
$ws.Range("I136").Value = 27570226
$ws.Range("J136").Value = 1779.6
$ws.Range("K136").Value = 82710678
$ws.Range("L136").Value = 5338.799999999999
$ws.Range("M136").Value = -82708128
$ws.Range("N136").Value = -10438.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 19049.908
$ws.Range("I5").Value = 36113.75
$ws.Range("J5").Value = 673.46155
$ws.Range("K5").Value = 108341.25
$ws.Range("L5").Value = 2020.38465
$ws.Range("M5").Value = -108229.25
$ws.Range("N5").Value = -2244.38465
$ws.Range("H39").Value = 2499.8
$ws.Range("J39").Value = 2499.8
$ws.Range("L39").Value = 7499.400000000001
$ws.Range("N39").Value = -8087.400000000001
$ws.Range("H55").Value = 2909
$ws.Range("J55").Value = 2909
$ws.Range("L55").Value = 8727
$ws.Range("N55").Value = -9081
$ws.Range("H107").Value = 1129.5508
$ws.Range("J107").Value = 2063.2144
$ws.Range("L107").Value = 6189.6432
$ws.Range("N107").Value = -10029.6432
$ws.Range("H131").Value = 564.8
$ws.Range("I131").Value = 522.1111
$ws.Range("K131").Value = 1566.3333
$ws.Range("M131").Value = 3473.6667
$ws.Range("H135").Value = 19049.908
$ws.Range("I135").Value = 36113.75
$ws.Range("J135").Value = 673.46155
$ws.Range("K135").Value = 325023.75
$ws.Range("L135").Value = 6061.15395
$ws.Range("M135").Value = -322488.75
$ws.Range("N135").Value = -11131.15395
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1044.3158
$ws.Range("I107").Value = 672.3077
$ws.Range("J107").Value = 1850.3334
$ws.Range("K107").Value = 672.3077
$ws.Range("L107").Value = 1850.3334
$ws.Range("N107").Value = -5690.3334
$ws.Range("H126").Value = 2780
$ws.Range("I126").Value = 1033.3334
$ws.Range("K126").Value = 3100.0002
$ws.Range("M126").Value = -630.0001999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 26651.098
$ws.Range("I132").Value = 1530.3043
$ws.Range("J132").Value = 58749.89
$ws.Range("K132").Value = 4590.9129
$ws.Range("L132").Value = 176249.67
$ws.Range("M132").Value = -2060.9129
$ws.Range("N132").Value = -181309.67
$ws.Range("H136").Value = 54503.965
$ws.Range("I136").Value = 25687.238
$ws.Range("J136").Value = 147604.16
$ws.Range("K136").Value = 77061.71400000001
$ws.Range("L136").Value = 442812.48
$ws.Range("M136").Value = -74511.71400000001
$ws.Range("N136").Value = -447912.48
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 23655.4
$ws.Range("J109").Value = 23655.4
$ws.Range("L109").Value = 23655.4
$ws.Range("N109").Value = -26429.4
$ws.Range("H132").Value = 64560.125
$ws.Range("I132").Value = 40454.46
$ws.Range("J132").Value = 169018
$ws.Range("K132").Value = 121363.38
$ws.Range("L132").Value = 507054
$ws.Range("M132").Value = -118833.38
$ws.Range("N132").Value = -512114
$ws.Range("H136").Value = 48491.934
$ws.Range("I136").Value = 30151
$ws.Range("J136").Value = 131026.125
$ws.Range("K136").Value = 90453
$ws.Range("L136").Value = 393078.375
$ws.Range("M136").Value = -87903
$ws.Range("N136").Value = -398178.375
